# edit.ps1
# Reproduce the OOXML diff for jira/exmple/example.xlsx:
#  - The old layout had the date (B1:D1, merged) above the TIME/TASK/SUMMARY/JIRA
#    header labels (row 2). The new layout swaps that: header labels now sit on
#    row 1 (A1:D1) and the date block moves to row 2 (A2:D2, merged).
#  - Two new work-log rows were recorded for day 42803 (row 10: 20 min, task
#    AW-13 "dfsfsd") and a whole new day 42804 was appended with its own date
#    header (row 12, merged A12:D12) and a work-log row (row 13: 22 min, task
#    jira-3 "sdsfsd").
#  - Active selection ends on D13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlGeneral = 1

# ---------------------------------------------------------------------------
# 1) Capture the old header-label row (row 2: TIME/TASK/SUMMARY/JIRA) values
#    before anything is overwritten.
# ---------------------------------------------------------------------------
$timeLabel    = $ws.Range("A2").Value()
$taskLabel    = $ws.Range("B2").Value()
$summaryLabel = $ws.Range("C2").Value()
$jiraLabel    = $ws.Range("D2").Value()

# Old date value that lived in B1 (merged B1:D1).
$firstDate = $ws.Range("B1").Value()

# ---------------------------------------------------------------------------
# 2) Remove the old merge so we can freely rewrite row 1 / row 2.
# ---------------------------------------------------------------------------
$ws.Range("B1:D1").UnMerge()

# ---------------------------------------------------------------------------
# 3) Row 1 becomes the bold header-label row (was row 2).
# ---------------------------------------------------------------------------
$headerRow = $ws.Range("A1:D1")
$headerRow.Font.Name = "Times New Roman"
$headerRow.Font.Bold = $true
$headerRow.Font.Color = 255
$headerRow.Font.Size = 10
$headerRow.NumberFormat = "General"
$headerRow.HorizontalAlignment = $xlGeneral

$ws.Range("A1").Value = $timeLabel
$ws.Range("B1").Value = $taskLabel
$ws.Range("C1").Value = $summaryLabel
$ws.Range("D1").Value = $jiraLabel

# ---------------------------------------------------------------------------
# 4) Row 2 becomes the merged date row for the first day (was row 1, col B).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = $firstDate
$ws.Range("A2").NumberFormat = "m/d/yyyy"
$ws.Range("A2").HorizontalAlignment = $xlCenter
$ws.Range("A2").Font.Name = "Arial"
$ws.Range("A2").Font.Bold = $false
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("A2:D2").Merge()

# ---------------------------------------------------------------------------
# 5) Row 10: existing work-log entry gains a duration + task reference.
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = 20
$ws.Range("C10").Value = "dfsfsd"
$ws.Range("D10").Value = "AW-13"

# ---------------------------------------------------------------------------
# 6) Row 12: brand-new merged date header for the second day (42804).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 42804
$ws.Range("A12").NumberFormat = "m/d/yyyy"
$ws.Range("A12").HorizontalAlignment = $xlCenter
$ws.Range("A12").Font.Name = "Arial"
$ws.Range("A12").Font.Bold = $false
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("A12:D12").Merge()

# ---------------------------------------------------------------------------
# 7) Row 13: new work-log entry for the second day.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 22
$ws.Range("C13").Value = "sdsfsd"
$ws.Range("D13").Value = "jira-3"

# ---------------------------------------------------------------------------
# 8) Restore the on-screen selection to D13, matching the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("D13").Select()
